# Update latest output (run 78)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Schedule": a new pump-run row is inserted between the existing two
# rows, shifting the old row 3 down to row 4 and updating values throughout.
# ---------------------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

# Insert a new row above the current row 3, pushing the existing row 3 (and
# its formatting) down to row 4.
$schedule.Rows.Item(3).Insert()

# Row 2 (updated values)
$schedule.Cells.Item(2, 1).Value = 46040.27083333334
$schedule.Cells.Item(2, 2).Value = 46040.4375
$schedule.Cells.Item(2, 3).Value = 4
$schedule.Cells.Item(2, 4).Value = 15.12
$schedule.Cells.Item(2, 5).Value = 140.93824875
$schedule.Cells.Item(2, 6).Value = 9.321312748015874

# Row 3 (brand new row)
$schedule.Cells.Item(3, 1).Value = 46040.45833333334
$schedule.Cells.Item(3, 2).Value = 46040.79166666666
$schedule.Cells.Item(3, 3).Value = 8
$schedule.Cells.Item(3, 4).Value = 30.24
$schedule.Cells.Item(3, 5).Value = -13.20846150000001
$schedule.Cells.Item(3, 6).Value = -0.4367877480158732

# Row 4 (previously row 3 - values are unchanged except Cost / Unit Cost)
$schedule.Cells.Item(4, 1).Value = 46040.83333333334
$schedule.Cells.Item(4, 2).Value = 46041
$schedule.Cells.Item(4, 3).Value = 4
$schedule.Cells.Item(4, 4).Value = 15.12
$schedule.Cells.Item(4, 5).Value = 333.5945925
$schedule.Cells.Item(4, 6).Value = 22.06313442460317

# ---------------------------------------------------------------------------
# Sheet "Detailed": refreshed price / status data for the latest run.
# ---------------------------------------------------------------------------
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Cells.Item(15, 5).Value = "ON"

$detailed.Cells.Item(16, 2).Value = 24.39469

$detailed.Cells.Item(18, 2).Value = -5.51

$detailed.Cells.Item(19, 2).Value = 36.06
$detailed.Cells.Item(19, 3).Value = "historical"

$detailed.Cells.Item(20, 2).Value = 20.88924
$detailed.Cells.Item(20, 3).Value = "historical"

$detailed.Cells.Item(21, 2).Value = 20.4213

$detailed.Cells.Item(22, 2).Value = 11.72682

$detailed.Cells.Item(23, 2).Value = 36.06
$detailed.Cells.Item(23, 5).Value = "OFF"

$detailed.Cells.Item(24, 2).Value = 0.7

$detailed.Cells.Item(25, 2).Value = 0.02018

$detailed.Cells.Item(26, 2).Value = 0.7

$detailed.Cells.Item(27, 2).Value = 0.7

$detailed.Cells.Item(28, 2).Value = -0.9016999999999999

$detailed.Cells.Item(29, 2).Value = -0.8791099999999999

$detailed.Cells.Item(30, 2).Value = -5.51011

$detailed.Cells.Item(31, 2).Value = -5.1729

$detailed.Cells.Item(32, 2).Value = -6.13636

$detailed.Cells.Item(33, 2).Value = -7.49401

$detailed.Cells.Item(34, 2).Value = -6.81821

$detailed.Cells.Item(35, 2).Value = -6.44839

$detailed.Cells.Item(36, 2).Value = -6.4473

$detailed.Cells.Item(37, 2).Value = 4.93756

$detailed.Cells.Item(38, 2).Value = 6.89183

$detailed.Cells.Item(39, 2).Value = 18.31138

$detailed.Cells.Item(41, 2).Value = 36.90374

$detailed.Cells.Item(42, 2).Value = 46.53487

$detailed.Cells.Item(43, 2).Value = 37.61157

$detailed.Cells.Item(44, 2).Value = 56.94145

$detailed.Cells.Item(45, 2).Value = 36.2

$detailed.Cells.Item(46, 2).Value = 56.68021
